$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for column A (Attribute) and B (Type), rows 2-21,
# reflecting the reordering/renaming described by the diff.
$data = @(
    @("case:concept:name", "str"),
    @("requested_service_url", "str"),
    @("human_workstation_green_button_pressed", "float"),
    @("lifecycle:state", "str"),
    @("response_status_code", "float"),
    @("process_model_id", "str"),
    @("operation_end_time", "datetime"),
    @("time:timestamp", "datetime"),
    @("org:resource", "str"),
    @("planned_operation_time", "str"),
    @("parameters", "dict"),
    @("unsatisfied_condition_description", "str"),
    @("event_id", "str"),
    @("SubProcessID", "str"),
    @("complete_service_time", "str"),
    @("current_task", "str"),
    @("concept:name", "str"),
    @("case", "str"),
    @("lifecycle:transition", "str"),
    @("identifier:id", "str")
)

$row = 2
foreach ($pair in $data) {
    $ws.Cells.Item($row, 1).Value = $pair[0]
    $ws.Cells.Item($row, 2).Value = $pair[1]
    $row++
}
